# Update metrics values on the worksheet: every data row (2-26) gets the
# same new set of values in columns B through Q (matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = 0.9999883732771242
    "C" = 0.9990763349190175
    "D" = 0.9999972881314654
    "E" = 0.9999711201841239
    "F" = 0.9999769148147516
    "G" = 0.00001085303545036088
    "H" = 0.000862200808886027
    "I" = 0.000001109367682337309
    "J" = 0.00005679383571578636
    "K" = 0.00002895160169906183
    "L" = 0.0001805656912053873
    "M" = 0.00329439454989242
    "N" = 0.9999069862169936
    "O" = 0.003434643746218193
    "P" = 64.86213150418904
    "Q" = 90.45852382642124
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col" + "$row").Value = $newValues[$col]
    }
}
